# Swap the deck's applied design from the "Integral" theme back to the
# built-in "Office Theme" palette (Design > Variants > Colors > Office).
#
# The Integral and Office Theme built-ins share an identical font scheme
# and format scheme (fills/lines/effects/background styles) - only the
# twelve theme colors differ: dk1, lt1, dk2, lt2, accent1-6, hlink and
# folHlink. Re-pointing every theme color at the Office Theme's RGB
# values reproduces the applied-design change across every slide, since
# they all share the one slide master's theme.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456   # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797  # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477   # folHlink -> 954F72
